# Update the single data row of the "Effects" sheet: replace the generic
# "liver effect / test effect" placeholder row with a concrete
# AChE-inhibition effect, and mark the record as synthetic test data
# (commit: "add synthetic data set simple vs real").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Effects")

$ws.Range("C2").Value = "AChE inhibition"
$ws.Range("D2").Value = "Acetylcholinesterase inhibition in liver tissue"
$ws.Range("F2").Value = "Neurotransmission"
$ws.Range("G2").Value = "Acetylcholinesterase"
$ws.Range("H2").Value = "Inhibition"
$ws.Range("I2").Value = "Liver"
$ws.Range("M2").Value = "Synthetic test data"
$ws.Range("P2").Value = $true
